$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 22, shifting existing rows 22-132 down to 23-133.
$ws.Rows.Item(22).Insert()

# Populate the new row 22 with the new weekly record.
$ws.Cells.Item(22, 1).Value = 4
$ws.Cells.Item(22, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(22, 3).Value = "Los Lagos"
$ws.Cells.Item(22, 4).Value = 44503
$ws.Cells.Item(22, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(22, 5).Value = 10
$ws.Cells.Item(22, 6).Value = 100112028
$ws.Cells.Item(22, 7).Value = "Sandia"
$ws.Cells.Item(22, 8).Value = "Sin especificar"
$ws.Cells.Item(22, 9).Value = "Primera"
$ws.Cells.Item(22, 10).Value = 150
$ws.Cells.Item(22, 11).Value = 1000
$ws.Cells.Item(22, 12).Value = 1000
$ws.Cells.Item(22, 13).Value = 1000
$ws.Cells.Item(22, 14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(22, 15).Value = "Perú"
$ws.Cells.Item(22, 16).Value = 1000
$ws.Cells.Item(22, 17).Value = 1
$ws.Cells.Item(22, 18).Value = "Hortaliza"
